$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1006.8571
$ws.Range("I98").Value = 1040.4706
$ws.Range("J98").Value = 864
$ws.Range("K98").Value = 1040.4706
$ws.Range("L98").Value = 864
$ws.Range("M98").Value = 457.5293999999999
$ws.Range("N98").Value = -3860

$ws.Range("H122").Value = 1006.8571
$ws.Range("I122").Value = 1040.4706
$ws.Range("J122").Value = 864
$ws.Range("K122").Value = 3121.4118
$ws.Range("L122").Value = 2592
$ws.Range("M122").Value = -671.4118000000003
$ws.Range("N122").Value = -7492

$ws.Range("H125").Value = 3472.6667
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 3472.6667
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 31254.0003
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -36174.0003

$ws.Range("H132").Value = 27780614
$ws.Range("I132").Value = 3002.25
$ws.Range("J132").Value = 41669420
$ws.Range("K132").Value = 9006.75
$ws.Range("L132").Value = 125008260
$ws.Range("M132").Value = -6476.75
$ws.Range("N132").Value = -125013320

$ws.Range("H136").Value = 24351.428
$ws.Range("J136").Value = 24351.428
$ws.Range("L136").Value = 24351.428
$ws.Range("N136").Value = -34551.428

$ws.Range("H137").Value = 1424.1111
$ws.Range("I137").Value = 1359.5385
$ws.Range("K137").Value = 4078.6155
$ws.Range("M137").Value = -1528.6155

$ws.Range("H138").Value = 2854.8767
$ws.Range("I138").Value = 761.3570999999999
$ws.Range("J138").Value = 5691.2583
$ws.Range("K138").Value = 2284.0713
$ws.Range("L138").Value = 17073.7749
$ws.Range("M138").Value = 2855.9287
$ws.Range("N138").Value = -27353.7749

$ws.Range("H139").Value = 29000.334
$ws.Range("J139").Value = 29000.334
$ws.Range("L139").Value = 29000.334
$ws.Range("N139").Value = -39280.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8584
$ws.Range("I45").Value = 10850.7
$ws.Range("J45").Value = 1028.3334
$ws.Range("K45").Value = 10850.7
$ws.Range("L45").Value = 1028.3334
$ws.Range("M45").Value = -10473.7
$ws.Range("N45").Value = -1782.3334

$ws.Range("H74").Value = 16667910
$ws.Range("I74").Value = 1213.1364
$ws.Range("K74").Value = 1213.1364
$ws.Range("M74").Value = -339.1364000000001

$ws.Range("H77").Value = 16667910
$ws.Range("I77").Value = 1213.1364
$ws.Range("K77").Value = 6065.682000000001
$ws.Range("M77").Value = -1697.682000000001

$ws.Range("H121").Value = 40000
$ws.Range("J121").Value = 40000
$ws.Range("L121").Value = 40000
$ws.Range("N121").Value = -43494

$ws.Range("H132").Value = 2510.561
$ws.Range("I132").Value = 1618
$ws.Range("J132").Value = 4667.5835
$ws.Range("K132").Value = 4854
$ws.Range("L132").Value = 14002.7505
$ws.Range("M132").Value = -2324
$ws.Range("N132").Value = -19062.7505

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 36921.668
$ws.Range("J53").Value = 36921.668
$ws.Range("L53").Value = 36921.668
$ws.Range("N53").Value = -38069.668

$ws.Range("H105").Value = 2136.4285
$ws.Range("I105").Value = 2136.4285
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2136.4285
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -389.4285
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4274848.5
$ws.Range("I16").Value = 6994006.5
$ws.Range("J16").Value = 1885.7142
$ws.Range("K16").Value = 6994006.5
$ws.Range("L16").Value = 1885.7142
$ws.Range("M16").Value = -6993719.5
$ws.Range("N16").Value = -2459.7142

$ws.Range("H58").Value = 4387240.5
$ws.Range("I58").Value = 6945355
$ws.Range("J58").Value = 1901.4286
$ws.Range("K58").Value = 6945355
$ws.Range("L58").Value = 1901.4286
$ws.Range("M58").Value = -6945152
$ws.Range("N58").Value = -2307.4286

$ws.Range("H86").Value = 2059.6
$ws.Range("I86").Value = 1600
$ws.Range("J86").Value = 2174.5
$ws.Range("K86").Value = 1600
$ws.Range("L86").Value = 2174.5
$ws.Range("M86").Value = -477
$ws.Range("N86").Value = -4420.5

$ws.Range("H89").Value = 2059.6
$ws.Range("I89").Value = 1600
$ws.Range("J89").Value = 2174.5
$ws.Range("K89").Value = 8000
$ws.Range("L89").Value = 10872.5
$ws.Range("M89").Value = -2384
$ws.Range("N89").Value = -22104.5

$ws.Range("H99").Value = 6949195.5
$ws.Range("I99").Value = 3470
$ws.Range("K99").Value = 3470
$ws.Range("M99").Value = -1972

$ws.Range("H113").Value = 4274848.5
$ws.Range("I113").Value = 6994006.5
$ws.Range("J113").Value = 1885.7142
$ws.Range("K113").Value = 6994006.5
$ws.Range("L113").Value = 1885.7142
$ws.Range("M113").Value = -6991836.5
$ws.Range("N113").Value = -6225.7142

$ws.Range("H122").Value = 1212.1818
$ws.Range("I122").Value = 1122.6666
$ws.Range("J122").Value = 1319.6
$ws.Range("K122").Value = 3367.9998
$ws.Range("L122").Value = 3958.8
$ws.Range("M122").Value = -917.9998000000001
$ws.Range("N122").Value = -8858.799999999999

$ws.Range("H126").Value = 6949195.5
$ws.Range("I126").Value = 3470
$ws.Range("K126").Value = 10410
$ws.Range("M126").Value = -7940

$ws.Range("H136").Value = 4387240.5
$ws.Range("I136").Value = 6945355
$ws.Range("J136").Value = 1901.4286
$ws.Range("K136").Value = 20836065
$ws.Range("L136").Value = 5704.2858
$ws.Range("M136").Value = -20833515
$ws.Range("N136").Value = -10804.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4846583.5
$ws.Range("I5").Value = 768.3
$ws.Range("J5").Value = 9692399
$ws.Range("K5").Value = 2304.9
$ws.Range("L5").Value = 29077197
$ws.Range("M5").Value = -2192.9
$ws.Range("N5").Value = -29077421

$ws.Range("H50").Value = 1306.6666
$ws.Range("I50").Value = 20
$ws.Range("J50").Value = 1950
$ws.Range("K50").Value = 60
$ws.Range("L50").Value = 5850
$ws.Range("M50").Value = 421
$ws.Range("N50").Value = -6812

$ws.Range("H53").Value = 1306.6666
$ws.Range("I53").Value = 20
$ws.Range("J53").Value = 1950
$ws.Range("K53").Value = 60
$ws.Range("L53").Value = 5850
$ws.Range("M53").Value = 421
$ws.Range("N53").Value = -6812

$ws.Range("H122").Value = 593.5454999999999
$ws.Range("I122").Value = 341.25
$ws.Range("J122").Value = 1266.3334
$ws.Range("K122").Value = 3071.25
$ws.Range("L122").Value = 11397.0006
$ws.Range("M122").Value = -621.25
$ws.Range("N122").Value = -16297.0006

$ws.Range("H125").Value = 3607.4443
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 3607.4443
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 10822.3329
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -20662.3329

$ws.Range("H132").Value = 13890298
$ws.Range("I132").Value = 1103.6
$ws.Range("J132").Value = 37038950
$ws.Range("K132").Value = 9932.4
$ws.Range("L132").Value = 333350550
$ws.Range("M132").Value = -7402.4
$ws.Range("N132").Value = -333355610

$ws.Range("H133").Value = 42311.414
$ws.Range("I133").Value = 147004.42
$ws.Range("J133").Value = 9000
$ws.Range("K133").Value = 441013.26
$ws.Range("L133").Value = 27000
$ws.Range("M133").Value = -435953.26
$ws.Range("N133").Value = -37120

$ws.Range("H134").Value = 10606.862
$ws.Range("I134").Value = 11863.637
$ws.Range("K134").Value = 35590.911
$ws.Range("M134").Value = -30520.911

$ws.Range("H135").Value = 4846583.5
$ws.Range("I135").Value = 768.3
$ws.Range("J135").Value = 9692399
$ws.Range("K135").Value = 6914.7
$ws.Range("L135").Value = 87231591
$ws.Range("M135").Value = -4379.7
$ws.Range("N135").Value = -87236661

$ws.Range("H137").Value = 12355778
$ws.Range("I137").Value = 12865.556
$ws.Range("J137").Value = 18527234
$ws.Range("K137").Value = 38596.66800000001
$ws.Range("L137").Value = 55581702
$ws.Range("M137").Value = -33496.66800000001
$ws.Range("N137").Value = -55591902

$ws.Range("H139").Value = 4803.227
$ws.Range("I139").Value = 11366
$ws.Range("K139").Value = 34098
$ws.Range("M139").Value = -28958

$ws.Range("H140").Value = 2065.3157
$ws.Range("I140").Value = 2068.9443
$ws.Range("K140").Value = 6206.8329
$ws.Range("M140").Value = -1026.8329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 8775286
$ws.Range("I132").Value = 15154781
$ws.Range("J132").Value = 3480
$ws.Range("K132").Value = 45464343
$ws.Range("L132").Value = 10440
$ws.Range("M132").Value = -45461813
$ws.Range("N132").Value = -15500

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 587.375
$ws.Range("I16").Value = 607.53845
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 607.53845
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -437.53845
$ws.Range("N16").Value = -840

$ws.Range("H61").Value = 3308.375
$ws.Range("I61").Value = 2605.5
$ws.Range("K61").Value = 2605.5
$ws.Range("M61").Value = -2403.5

$ws.Range("H113").Value = 3308.375
$ws.Range("I113").Value = 2605.5
$ws.Range("K113").Value = 2605.5
$ws.Range("M113").Value = -435.5
